$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Valor Mora" total (row 11) ---
$ws.Range("E11").Value = 173333

# --- Update "Cant. Periodos" (row 13) : 2 -> 3 periods ---
$ws.Range("F13").Value = 3

# --- Insert a new detail row (18) below the existing two detail rows (16,17) so
#     a third debt period can be listed. Formatting now needs to cascade down one
#     band: the new row 18 takes the old row 17 look, and row 17 takes row 16's. ---
$ws.Rows.Item(18).Insert()

# Row 18 inherits the (pre-shift) row 17 formatting
$ws.Range("B17:J17").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 17 now inherits row 16's formatting
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 16: newest period (2507), value now 80000 ---
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 80000

# Row 17 keeps the prior 2506/80000 entry (value untouched, only format changed above).

# --- Row 18 (new): oldest period (2505), value 13333 ---
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "92515451"
$ws.Range("D18").Value = "LUIS ALONSO MERCADO PEÑATE"
$ws.Range("E18").Value = "2505"
$ws.Range("F18").Value = 13333
$ws.Range("G18").Value = 2000000
